$wb = $excel.ActiveWorkbook

function Get-HyperlinkMap($ws) {
    $map = @{}
    foreach ($h in $ws.Hyperlinks) {
        $addr = $h.Range.Address()
        $map[$addr] = @{ Address = $h.Address; Display = $h.TextToDisplay }
    }
    return $map
}

# ---------------------------------------------------------------------------
# 1) Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This shared string is used on the Overview sheet (B2,C2,B3,C3) and on
#    the per-language sheets (C2,C3).
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: handback completed - the target file now also counts as
#    the handback file. Populate F (Latest Target File) / G (Latest Handback
#    File) with hyperlinks mirroring the existing A (Source File Name) /
#    D (Latest Handoff File) hyperlinks, and refresh the handback timestamp.
# ---------------------------------------------------------------------------
$zhLinks = Get-HyperlinkMap $wsZhCn

$zhA2 = $zhLinks["`$A`$2"]
$zhD2 = $zhLinks["`$D`$2"]
$zhA3 = $zhLinks["`$A`$3"]
$zhD3 = $zhLinks["`$D`$3"]

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), $zhA2.Address, "", "", $zhA2.Display) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhD2.Address, "", "", $zhD2.Display) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), $zhA3.Address, "", "", $zhA3.Display) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhD3.Address, "", "", $zhD3.Display) | Out-Null

$wsZhCn.Range("F2:G3").Font.Underline = 2
$wsZhCn.Range("F2:G3").Font.Color = 15570276
$wsZhCn.Range("F2:G3").Font.Name = "Calibri"
$wsZhCn.Range("F2:G3").Font.Size = 11

$wsZhCn.Range("H2").Value = "2016-03-18 03:44:25"
$wsZhCn.Range("H3").Value = "2016-03-18 03:44:25"

# ---------------------------------------------------------------------------
# 3) de-de sheet: same handback treatment, but completed a little later.
# ---------------------------------------------------------------------------
$deLinks = Get-HyperlinkMap $wsDeDe

$deA2 = $deLinks["`$A`$2"]
$deD2 = $deLinks["`$D`$2"]
$deA3 = $deLinks["`$A`$3"]
$deD3 = $deLinks["`$D`$3"]

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), $deA2.Address, "", "", $deA2.Display) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deD2.Address, "", "", $deD2.Display) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), $deA3.Address, "", "", $deA3.Display) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deD3.Address, "", "", $deD3.Display) | Out-Null

$wsDeDe.Range("F2:G3").Font.Underline = 2
$wsDeDe.Range("F2:G3").Font.Color = 15570276
$wsDeDe.Range("F2:G3").Font.Name = "Calibri"
$wsDeDe.Range("F2:G3").Font.Size = 11

$wsDeDe.Range("H2").Value = "2016-03-18 03:44:39"
$wsDeDe.Range("H3").Value = "2016-03-18 03:44:39"

$wb.Save()
